# fix validate excel + add payment success
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data: FlightNumber changes from QH3456 -> QH4444 (AirplaneCode QH1111 stays the same)
$ws.Range("A2").Value = "QH4444"
$ws.Range("B2").Value = "QH1111"

# DepartureTime shifts from 11/6/2024 08:00 to 11/6/2024 20:00 (+12h / +0.5 day)
$ws.Range("C2").Value = 45602.833333333336

# Update the active selection to E8 (was C3)
[void]$ws.Range("E8").Select()
